$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new task ("Priorité 1" list) in row 6, following the same layout
# as the existing rows (merged A:D, centered + wrapped text).
$ws.Range("A6").Value = "Ranger le code (supprimer des classes si pas nécessaires)"
$ws.Range("A6:D6").Merge()

# Match the row height used for this wrapped entry.
$ws.Rows.Item(6).RowHeight = 30.75

# Leave the selection where the author ended up after the edit.
$ws.Range("G7").Select()
